$wb = $excel.ActiveWorkbook

$wsA = $wb.Worksheets.Item("展览")
$wsB = $wb.Worksheets.Item("演出")
$wsD = $wb.Worksheets.Item("全部类型")

$wsA.Range("F12").Value = 628
$wsA.Range("F14").Value = 527
$wsA.Range("F15").Value = 378
$wsA.Range("F18").Value = 1330
$wsA.Range("F20").Value = 1603
$wsA.Range("F25").Value = 532
$wsA.Range("F28").Value = 99
$wsA.Range("G30").Value = 45
$wsA.Range("F32").Value = 3815
$wsA.Range("F34").Value = 69
$wsA.Range("F35").Value = 781
$wsA.Range("F37").Value = 1814
$wsB.Range("F3").Value = 42
$wsD.Range("F12").Value = 628
$wsD.Range("F14").Value = 527
$wsD.Range("F16").Value = 378
$wsD.Range("F17").Value = 134
$wsD.Range("F19").Value = 1330
$wsD.Range("F21").Value = 1603
$wsD.Range("F26").Value = 532
$wsD.Range("F29").Value = 99
$wsD.Range("G31").Value = 45
$wsD.Range("F33").Value = 3815
$wsD.Range("F34").Value = 42
$wsD.Range("F36").Value = 69
$wsD.Range("F37").Value = 781
$wsD.Range("F39").Value = 1814
